{"js": "// Update the date heading paragraph, then the 20x5 grid of arithmetic\n// problems inside the single table, cell by cell in row-major order.\nconst NEW_DATE = \"2024-02-23 Friday\";\nconst NEW_CELLS = [\"16+77=\", \"91-53=\", \"33-28=\", \"17+37=\", \"90-57=\", \"40-28=\", \"31-16=\", \"55+27=\", \"61-28=\", \"16+28=\", \"57+14=\", \"48+3=\", \"84-75=\", \"27+18=\", \"90-24=\", \"19+19=\", \"52-18=\", \"38+28=\", \"39+5=\", \"8+14=\", \"19+55=\", \"28+34=\", \"74-36=\", \"37+19=\", \"70-22=\", \"34+48=\", \"53-5=\", \"8+4=\", \"47+38=\", \"5+79=\", \"63-36=\", \"73-56=\", \"14+19=\", \"43-28=\", \"7+29=\", \"43-14=\", \"28+66=\", \"61-53=\", \"40-13=\", \"90-15=\", \"18+15=\", \"27+17=\", \"28+65=\", \"60-2=\", \"81-46=\", \"26+56=\", \"61-15=\", \"64+28=\", \"65-17=\", \"6+45=\", \"9+82=\", \"28+33=\", \"6+35=\", \"46+7=\", \"31-15=\", \"52-3=\", \"45+16=\", \"9+25=\", \"65-7=\", \"41-22=\", \"37+56=\", \"77+4=\", \"18+48=\", \"29+7=\", \"18+73=\", \"59+16=\", \"8+65=\", \"34-5=\", \"45+8=\", \"61-36=\", \"77+19=\", \"74-29=\", \"7+68=\", \"9+27=\", \"72-15=\", \"60-1=\", \"91-2=\", \"14+8=\", \"54-8=\", \"94-36=\", \"49+9=\", \"37+45=\", \"40-39=\", \"50-33=\", \"7+79=\", \"37+57=\", \"40-2=\", \"4+59=\", \"35+18=\", \"39+34=\", \"28+66=\", \"46+26=\", \"43-17=\", \"19+22=\", \"56+6=\", \"45+28=\", \"52-33=\", \"80-56=\", \"8+43=\", \"90-51=\"];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load('items');\nconst tables = body.tables;\ntables.load('items');\nawait context.sync();\n\n// First paragraph holds the date string, e.g. '2024-02-22 Thursday'.\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load('text');\nawait context.sync();\ndateParagraph.insertText(NEW_DATE, Word.InsertLocation.replace);\n\n// The single table holds the math problems, 20 rows x 5 columns.\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load('items');\nawait context.sync();\n\nlet cellIdx = 0;\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load('items');\n  await context.sync();\n  for (const cell of cells.items) {\n    cell.value = NEW_CELLS[cellIdx];\n    cellIdx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading paragraph, then the 20x5 grid of arithmetic\n# problems inside the single table, cell by cell in row-major order.\n$d = $word.ActiveDocument\n\n$newDate = \"2024-02-23 Friday\"\n$d.Paragraphs.Item(1).Range.Text = $newDate\n\n$newCells = @(\n    \"16+77=\", \"91-53=\", \"33-28=\", \"17+37=\", \"90-57=\",\n    \"40-28=\", \"31-16=\", \"55+27=\", \"61-28=\", \"16+28=\",\n    \"57+14=\", \"48+3=\", \"84-75=\", \"27+18=\", \"90-24=\",\n    \"19+19=\", \"52-18=\", \"38+28=\", \"39+5=\", \"8+14=\",\n    \"19+55=\", \"28+34=\", \"74-36=\", \"37+19=\", \"70-22=\",\n    \"34+48=\", \"53-5=\", \"8+4=\", \"47+38=\", \"5+79=\",\n    \"63-36=\", \"73-56=\", \"14+19=\", \"43-28=\", \"7+29=\",\n    \"43-14=\", \"28+66=\", \"61-53=\", \"40-13=\", \"90-15=\",\n    \"18+15=\", \"27+17=\", \"28+65=\", \"60-2=\", \"81-46=\",\n    \"26+56=\", \"61-15=\", \"64+28=\", \"65-17=\", \"6+45=\",\n    \"9+82=\", \"28+33=\", \"6+35=\", \"46+7=\", \"31-15=\",\n    \"52-3=\", \"45+16=\", \"9+25=\", \"65-7=\", \"41-22=\",\n    \"37+56=\", \"77+4=\", \"18+48=\", \"29+7=\", \"18+73=\",\n    \"59+16=\", \"8+65=\", \"34-5=\", \"45+8=\", \"61-36=\",\n    \"77+19=\", \"74-29=\", \"7+68=\", \"9+27=\", \"72-15=\",\n    \"60-1=\", \"91-2=\", \"14+8=\", \"54-8=\", \"94-36=\",\n    \"49+9=\", \"37+45=\", \"40-39=\", \"50-33=\", \"7+79=\",\n    \"37+57=\", \"40-2=\", \"4+59=\", \"35+18=\", \"39+34=\",\n    \"28+66=\", \"46+26=\", \"43-17=\", \"19+22=\", \"56+6=\",\n    \"45+28=\", \"52-33=\", \"80-56=\", \"8+43=\", \"90-51=\"\n)\n\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $t.Cell($r, $c).Range.Text = $newCells[$idx]\n        $idx = $idx + 1\n    }\n}\n"}
